# Auto-generated edit script
# Updates the "想去人数" (F) and "最低票价" (G) figures across the
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets,
# matching the refreshed scrape output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 235
$ws1.Range("F5").Value = 2855
$ws1.Range("F8").Value = 2206
$ws1.Range("F9").Value = 337
$ws1.Range("F13").Value = 2540
$ws1.Range("F15").Value = 1322
$ws1.Range("F16").Value = 4664
$ws1.Range("F18").Value = 5037
$ws1.Range("F19").Value = 1588
$ws1.Range("F20").Value = 2847
$ws1.Range("G20").Value = 6.6
$ws1.Range("F21").Value = 3241
$ws1.Range("F22").Value = 160
$ws1.Range("F23").Value = 1538
$ws1.Range("F24").Value = 250
$ws1.Range("F25").Value = 835
$ws1.Range("F26").Value = 102
$ws1.Range("F27").Value = 285
$ws1.Range("F28").Value = 961
$ws1.Range("F29").Value = 1785
$ws1.Range("F30").Value = 115
$ws1.Range("F31").Value = 274
$ws1.Range("F32").Value = 679
$ws1.Range("F33").Value = 155
$ws1.Range("F34").Value = 326
$ws1.Range("F35").Value = 400

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 89

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 235
$ws4.Range("F11").Value = 2855
$ws4.Range("F13").Value = 2206
$ws4.Range("F14").Value = 337
$ws4.Range("F15").Value = 89
$ws4.Range("F21").Value = 2540
$ws4.Range("F22").Value = 1322
$ws4.Range("F26").Value = 4664
$ws4.Range("F28").Value = 5037
$ws4.Range("F29").Value = 1588
$ws4.Range("F30").Value = 2847
$ws4.Range("G30").Value = 6.6
$ws4.Range("F31").Value = 3241
$ws4.Range("F32").Value = 160
$ws4.Range("F35").Value = 1538
$ws4.Range("F37").Value = 250
$ws4.Range("F38").Value = 835
$ws4.Range("F39").Value = 102
$ws4.Range("F40").Value = 285
$ws4.Range("F41").Value = 961
$ws4.Range("F43").Value = 1785
$ws4.Range("F44").Value = 115
$ws4.Range("F45").Value = 274
$ws4.Range("F46").Value = 679
$ws4.Range("F47").Value = 155
$ws4.Range("F48").Value = 326
$ws4.Range("F49").Value = 400

